$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1424851084621821
$ws.Range("D2").Value = 0.02162275288967663
$ws.Range("E2").Value = 0.9706460231672907
$ws.Range("F2").Value = 0.2074658241646503
$ws.Range("G2").Value = 0.0959412586952908
$ws.Range("H2").Value = 0.2308347742849293
$ws.Range("M2").Value = 5.463493702426234
$ws.Range("O2").Value = 0.5408165625539283
$ws.Range("B3").Value = 0.1329187634610349
$ws.Range("D3").Value = 0.0190173867394634
$ws.Range("E3").Value = 0.8550687529744039
$ws.Range("F3").Value = 0.2048379074659934
$ws.Range("G3").Value = 0.09301119831525284
$ws.Range("H3").Value = 0.2355292435101219
$ws.Range("M3").Value = 4.771238460786236
$ws.Range("O3").Value = 0.5440645604220151
$ws.Range("B4").Value = 0.1271150842600832
$ws.Range("D4").Value = 0.01740928393535768
$ws.Range("E4").Value = 0.7845186426024782
$ws.Range("F4").Value = 0.2036525758176424
$ws.Range("G4").Value = 0.0915470519350734
$ws.Range("H4").Value = 0.2387424300278767
$ws.Range("M4").Value = 4.344752309698066
$ws.Range("O4").Value = 0.5474501212801073
$ws.Range("B5").Value = 0.1247677545229493
$ws.Range("D5").Value = 0.01675193829008492
$ws.Range("E5").Value = 0.7558687716256713
$ws.Range("F5").Value = 0.2032753348683798
$ws.Range("G5").Value = 0.09103240660711265
$ws.Range("H5").Value = 0.2401340853268792
$ws.Range("M5").Value = 4.170585153446837
$ws.Range("O5").Value = 0.5491729930788409
$ws.Range("B6").Value = 0.1243790550380339
$ws.Range("D6").Value = 0.01664266633542155
$ws.Range("E6").Value = 0.7511173939407882
$ws.Range("F6").Value = 0.2032190286948534
$ws.Range("G6").Value = 0.0909518385611392
$ws.Range("H6").Value = 0.24037010985856
$ws.Range("M6").Value = 4.14164222545989
$ws.Range("O6").Value = 0.5494796115546592
$ws.Range("B7").Value = 0.1270833554913935
$ws.Range("D7").Value = 0.01740042685526788
$ws.Range("E7").Value = 0.7841318615622583
$ws.Range("F7").Value = 0.2036470623590745
$ws.Range("G7").Value = 0.0915397821248547
$ws.Range("H7").Value = 0.2387608665370777
$ws.Range("M7").Value = 4.342404941363725
$ws.Range("O7").Value = 0.5474719755392101
$ws.Range("B8").Value = 0.1391721336851361
$ws.Range("D8").Value = 0.02072621232030514
$ws.Range("E8").Value = 0.9307060389488129
$ws.Range("F8").Value = 0.2064697025229165
$ws.Range("G8").Value = 0.09486012558110701
$ws.Range("H8").Value = 0.2323842252968547
$ws.Range("M8").Value = 5.225096595408729
$ws.Range("O8").Value = 0.5416437217200496
$ws.Range("B9").Value = 0.1634320469635782
$ws.Range("D9").Value = 0.02717833388996382
$ws.Range("E9").Value = 1.221667686668354
$ws.Range("F9").Value = 0.2154882714520951
$ws.Range("G9").Value = 0.1041289242423247
$ws.Range("H9").Value = 0.2225452438570699
$ws.Range("M9").Value = 6.945215016769282
$ws.Range("O9").Value = 0.5415525096821767
$ws.Range("B10").Value = 0.181592800107822
$ws.Range("D10").Value = 0.03187268654438924
$ws.Range("E10").Value = 1.437977172198828
$ws.Range("F10").Value = 0.2243572009104469
$ws.Range("G10").Value = 0.1127591527459657
$ws.Range("H10").Value = 0.2169982298749176
$ws.Range("M10").Value = 8.203440162718323
$ws.Range("O10").Value = 0.5488119871852177
$ws.Range("B11").Value = 0.1899278958854751
$ws.Range("D11").Value = 0.03399762527699579
$ws.Range("E11").Value = 1.53702969910816
$ws.Range("F11").Value = 0.2289049647439683
$ws.Range("G11").Value = 0.117110952040349
$ws.Range("H11").Value = 0.2148522142552167
$ws.Range("M11").Value = 8.774913760924619
$ws.Range("O11").Value = 0.5537957005988972
$ws.Range("B12").Value = 0.193094748267427
$ws.Range("D12").Value = 0.0348007081958599
$ws.Range("E12").Value = 1.574640335712957
$ws.Range("F12").Value = 0.2307030228395419
$ws.Range("G12").Value = 0.1188226151508758
$ws.Range("H12").Value = 0.2140948622500218
$ws.Range("M12").Value = 8.991209734362428
$ws.Range("O12").Value = 0.5559321587641932
$ws.Range("B13").Value = 0.1924122424264567
$ws.Range("D13").Value = 0.03462782148393728
$ws.Range("E13").Value = 1.566535562887566
$ws.Range("F13").Value = 0.2303123684148503
$ws.Range("G13").Value = 0.1184511033738858
$ws.Range("H13").Value = 0.2142554955739087
$ws.Range("M13").Value = 8.944631006893474
$ws.Range("O13").Value = 0.5554608297309471
$ws.Range("B14").Value = 0.1901882239814512
$ws.Range("D14").Value = 0.03406372754488984
$ws.Range("E14").Value = 1.540121865816076
$ws.Range("F14").Value = 0.2290513583419198
$ws.Range("G14").Value = 0.1172504799336025
$ws.Range("H14").Value = 0.2147887922355523
$ws.Range("M14").Value = 8.7927106124651
$ws.Range("O14").Value = 0.5539664302161214
$ws.Range("B15").Value = 0.188827318376255
$ws.Range("D15").Value = 0.03371799528760278
$ws.Range("E15").Value = 1.523956186395878
$ws.Range("F15").Value = 0.2282889032950877
$ws.Range("G15").Value = 0.1165234376194917
$ws.Range("H15").Value = 0.2151226841300371
$ws.Range("M15").Value = 8.699641482992831
$ws.Range("O15").Value = 0.5530837479842319
$ws.Range("B16").Value = 0.1810495591056309
$ws.Range("D16").Value = 0.03173359863458103
$ws.Range("E16").Value = 1.431517634398858
$ws.Range("F16").Value = 0.2240705018477556
$ws.Range("G16").Value = 0.1124835342275929
$ws.Range("H16").Value = 0.2171461585016914
$ws.Range("M16").Value = 8.16607663579066
$ws.Range("O16").Value = 0.5485207581914153
$ws.Range("B17").Value = 0.1762969799441834
$ws.Range("D17").Value = 0.03051348450337343
$ws.Range("E17").Value = 1.374981990498583
$ws.Range("F17").Value = 0.2216156002977172
$ws.Range("G17").Value = 0.1101160217309882
$ws.Range("H17").Value = 0.2184848756201347
$ws.Range("M17").Value = 7.838535522277766
$ws.Range("O17").Value = 0.5461574217926284
$ws.Range("B18").Value = 0.1735703603413299
$ws.Range("D18").Value = 0.02981071839761995
$ws.Range("E18").Value = 1.342525072306984
$ws.Range("F18").Value = 0.2202517869099907
$ws.Range("G18").Value = 0.1087942246885518
$ws.Range("H18").Value = 0.2192903256263037
$ws.Range("M18").Value = 7.650056512462072
$ws.Range("O18").Value = 0.5449558975923594
$ws.Range("B19").Value = 0.1726483674543715
$ws.Range("D19").Value = 0.02957260620625135
$ws.Range("E19").Value = 1.331545947629905
$ws.Range("F19").Value = 0.2197982342185725
$ws.Range("G19").Value = 0.1083534701274118
$ws.Range("H19").Value = 0.2195690952754177
$ws.Range("M19").Value = 7.586225391474215
$ws.Range("O19").Value = 0.5445759593917785
$ws.Range("B20").Value = 0.1768021825022998
$ws.Range("D20").Value = 0.03064347057387806
$ws.Range("E20").Value = 1.380993950000317
$ws.Range("F20").Value = 0.2218719269620451
$ws.Range("G20").Value = 0.1103638965609406
$ws.Range("H20").Value = 0.2183386899448436
$ws.Range("M20").Value = 7.873411593743072
$ws.Range("O20").Value = 0.54639261710048
$ws.Range("B21").Value = 0.1908411867247253
$ws.Range("D21").Value = 0.03422945916375397
$ws.Range("E21").Value = 1.547877386427274
$ws.Range("F21").Value = 0.2294196701795244
$ws.Range("G21").Value = 0.1176013824690614
$ws.Range("H21").Value = 0.2146306410099612
$ws.Range("M21").Value = 8.837336088608197
$ws.Range("O21").Value = 0.5543985475171098
$ws.Range("B22").Value = 0.2000778550140012
$ws.Range("D22").Value = 0.03656384036887061
$ws.Range("E22").Value = 1.657541594577197
$ws.Range("F22").Value = 0.2347960681961609
$ws.Range("G22").Value = 0.1227041526604609
$ws.Range("H22").Value = 0.212530021995974
$ws.Range("M22").Value = 9.466694220692375
$ws.Range("O22").Value = 0.5610871435861213
$ws.Range("B23").Value = 0.1951424749966577
$ws.Range("D23").Value = 0.03531880720652225
$ws.Range("E23").Value = 1.598954447615228
$ws.Range("F23").Value = 0.2318852967217069
$ws.Range("G23").Value = 0.1199457780204511
$ws.Range("H23").Value = 0.2136212800547099
$ws.Range("M23").Value = 9.130843465882094
$ws.Range("O23").Value = 0.5573815766325367
$ws.Range("B24").Value = 0.1765737626628692
$ws.Range("D24").Value = 0.03058470793408219
$ws.Range("E24").Value = 1.378275799862138
$ws.Range("F24").Value = 0.2217558937534889
$ws.Range("G24").Value = 0.1102517100737117
$ws.Range("H24").Value = 0.2184046690465919
$ws.Range("M24").Value = 7.857644656803245
$ws.Range("O24").Value = 0.5462857962556029
$ws.Range("B25").Value = 0.1568098952690775
$ws.Range("D25").Value = 0.02544070513221186
$ws.Range("E25").Value = 1.142538945851896
$ws.Range("F25").Value = 0.2126633081431208
$ws.Range("G25").Value = 0.1013119343309796
$ws.Range("H25").Value = 0.2249157493934034
$ws.Range("M25").Value = 6.480935901360112
$ws.Range("O25").Value = 0.5403207734865418

Write-Output "Updated 192 cells (B,D,E,F,G,H,M,O for rows 2-25) on Sheet1"
